# Updated symbol list on Mon Feb 13 19:56:53 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns for the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell stores a literal text string (not a number) in the source file,
# so values are entered apostrophe-prefixed (forces text) and the style is
# reset to "Normal" afterwards so no incidental number-format style gets
# attached to the cell (matches original, unstyled text cells).

$ws.Range("D2").Value = "'289.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-9.49%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.82%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.040"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-4.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07289"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-5.88%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.527"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-12.68%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9176"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.87%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.1175"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-5.56%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1719"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-8.81%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.08694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-4.96%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.04176"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.72%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'0.17%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.001274"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.71%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.005829"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.77%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.391"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.51%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'4.285"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.33%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'0.3319"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.19%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.888"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-9.91%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1353"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.03%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'0.03859"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.25%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'0.16%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.003853"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-6.53%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.86%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D38").Value = "'0.02313"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-9.82%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.04957"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.006558"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'229.48%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-1.14%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'-3.46%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'4.75%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007068"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-14.49%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3122"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.50%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006451"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.78%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-82.25%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.17%"
$ws.Range("E51").Style = "Normal"
